$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.500.97"
$ws.Range("E2").Value = "  +6.51%  "

$ws.Range("D3").Value = "3.349.31"
$ws.Range("E3").Value = "  +2.78%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.16"
$ws.Range("E5").Value = "  +4.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.35"
$ws.Range("E6").Value = "  +2.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.587"
$ws.Range("E7").Value = "  +4.76%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.638"
$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.72"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("E11").Value = "  +3.33%  "

$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").Value = "3.883.71"
$ws.Range("E13").Value = "  +2.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.52"
$ws.Range("E14").Value = "  +3.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.94"
$ws.Range("E15").Value = "  +4.40%  "

$ws.Range("D16").Value = "3.320.74"
$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "60.271.58"
$ws.Range("E18").Value = "  +6.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.84"
$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("E20").Value = "  +2.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000110"
$ws.Range("E21").Value = "  +4.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.25"
$ws.Range("E22").Value = "  +2.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "302.36"
$ws.Range("E23").Value = "  -1.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.76"
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.20"
$ws.Range("E25").Value = "  +1.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.68"
$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.77"
$ws.Range("E27").Value = "  +7.16%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.47"
$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.181"
$ws.Range("E29").Value = "  +6.18%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.97"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.116"
$ws.Range("E31").Value = "  +5.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.58"
$ws.Range("E32").Value = "  +20.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.53"
$ws.Range("E33").Value = "  +4.71%  "

$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.00"
$ws.Range("E35").Value = "  +6.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0508"
$ws.Range("E36").Value = "  +6.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.32"
$ws.Range("E37").Value = "  +1.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.10"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.41"
$ws.Range("E40").Value = "  -3.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "137.82"
$ws.Range("E41").Value = "  +1.94%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.296"
$ws.Range("E42").Value = "  +5.82%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.123"
$ws.Range("E43").Value = "  +2.48%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.93"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.97"
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.98"
$ws.Range("E46").Value = "  -1.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.28"
$ws.Range("E47").Value = "  +9.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.43"
$ws.Range("E48").Value = "  +1.64%  "

$ws.Range("D49").Value = "2.203.89"
$ws.Range("E49").Value = "  +2.65%  "

$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("E51").Value = "  -1.21%  "

